# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
# per scraped data refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'64.298.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.33%  '
$ws.Range("D3").Value = "'2.535.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.98%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'582.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").Value = "'152.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.56%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = "'0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("D11").Value = "'5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").Value = "'0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = "'29.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.92%  '
$ws.Range("D14").Value = "'0.0000180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.83%  '
$ws.Range("D15").Value = "'2.986.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.63%  '
$ws.Range("D16").Value = "'64.136.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.26%  '
$ws.Range("D17").Value = "'2.524.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").Value = "'7.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").Value = "'11.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = "'4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.45%  '
$ws.Range("D21").Value = "'328.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = "'2.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").Value = "'10.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").Value = "'65.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = "'662.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("D27").Value = "'0.0000103"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.35%  '
$ws.Range("D28").Value = "'2.646.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("E29").Value = '  +2.96%  '
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Value = "'8.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D33").Value = "'0.136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.25%  '
$ws.Range("D36").Value = "'4.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.80%  '
$ws.Range("D37").Value = "'5.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("D39").Value = "'18.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").Value = "'152.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("D41").Value = "'2.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.59%  '
$ws.Range("E42").Value = '  +3.43%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = "'158.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.97%  '
$ws.Range("D45").Value = "'0.0₆0303"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("D46").Value = "'15.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").Value = "'3.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.05%  '
$ws.Range("D48").Value = "'21.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.96%  '
$ws.Range("D49").Value = "'0.621"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.33%  '
$ws.Range("E50").Value = '  +2.07%  '
$ws.Range("E51").Value = '  +1.87%  '
